$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply "Text" number format first to cells whose new value would otherwise be
# auto-converted to a floating point number by Excel (so the literal text from
# the source data, e.g. "248.01" or "8.00", is preserved exactly as a string).

# Row 2
$ws.Range("D2").Value = "91.879.65"
$ws.Range("E2").Value = "  +0.85%  "

# Row 3
$ws.Range("D3").Value = "3.121.32"
$ws.Range("E3").Value = "  -0.72%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.01"
$ws.Range("E5").Value = "  +3.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.30"
$ws.Range("E6").Value = "  -2.12%  "

# Row 7
$ws.Range("E7").Value = "  +6.49%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.372"
$ws.Range("E8").Value = "  +1.62%  "

# Row 9
$ws.Range("E9").Value = "  +0.02%  "

# Row 10
$ws.Range("E10").Value = "  -0.75%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.759"
$ws.Range("E11").Value = "  +4.91%  "

# Row 12
$ws.Range("E12").Value = "  +3.03%  "

# Row 13
$ws.Range("E13").Value = "  +2.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.66"
$ws.Range("E14").Value = "  -3.71%  "

# Row 15
$ws.Range("B15").Value = "Toncoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.51"
$ws.Range("E15").Value = "  -0.68%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "91.670.23"
$ws.Range("E16").Value = "  +0.77%  "

# Row 17
$ws.Range("D17").Value = "3.699.83"
$ws.Range("E17").Value = "  -0.54%  "

# Row 18
$ws.Range("D18").Value = "3.114.92"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("E19").Value = "  -0.22%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.61"
$ws.Range("E20").Value = "  +1.91%  "

# Row 21
$ws.Range("E21").Value = "  +1.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.81"
$ws.Range("E22").Value = "  +2.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "448.11"
$ws.Range("E23").Value = "  +0.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.30"
$ws.Range("E24").Value = "  +3.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.92"
$ws.Range("E25").Value = "  -1.13%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.18"
$ws.Range("E26").Value = "  +1.11%  "

# Row 27
$ws.Range("E27").Value = "  -3.45%  "

# Row 28
$ws.Range("D28").Value = "3.258.19"
$ws.Range("E28").Value = "  -1.15%  "

# Row 29
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("E30").Value = "  +16.82%  "

# Row 31
$ws.Range("E31").Value = "  +21.56%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.37"
$ws.Range("E32").Value = "  -3.32%  "

# Row 33
$ws.Range("E33").Value = "  +16.58%  "

# Row 34
$ws.Range("E34").Value = "  +3.71%  "

# Row 35
$ws.Range("E35").Value = "  +30.50%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.00"
$ws.Range("E36").Value = "  +10.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.70"
$ws.Range("E37").Value = "  -1.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.18"
$ws.Range("E38").Value = "  +25.26%  "

# Row 39
$ws.Range("E39").Value = "  -0.63%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "498.14"
$ws.Range("E40").Value = "  -3.39%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("E41").Value = "  -3.79%  "

# Row 42
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.424"
$ws.Range("E43").Value = "  +1.51%  "

# Row 44
$ws.Range("E44").Value = "  -0.17%  "

# Row 46
$ws.Range("E46").Value = "  -0.72%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.701"
$ws.Range("E47").Value = "  +0.66%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "153.87"
$ws.Range("E48").Value = "  +1.70%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.55"
$ws.Range("E49").Value = "  -0.21%  "

# Row 50
$ws.Range("E50").Value = "  -0.85%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.42"
$ws.Range("E51").Value = "  -2.49%  "
